$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.75560066666667
$ws.Range("H2").Value = 50.266802
$ws.Range("I2").Value = 0.9064438825950116
$ws.Range("J2").Value = 0.9064438825950115
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.457778999999999
$ws.Range("N2").Value = 22.373337
$ws.Range("O2").Value = 0.08029647035915141
$ws.Range("P2").Value = 0.0802964703591514
$ws.Range("Q2").Value = 124.9595667842527
$ws.Range("R2").Value = 1124.636101058274
$ws.Range("S2").Value = 0.07278424435102447
$ws.Range("T2").Value = 0.07278424435102446

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.75560066666667
$ws.Range("H3").Value = 50.266802
$ws.Range("I3").Value = 0.9064438825950116
$ws.Range("J3").Value = 0.9064438825950115
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.732509
$ws.Range("N3").Value = 8.197526999999999
$ws.Range("O3").Value = 0.02942039820764526
$ws.Range("P3").Value = 0.02942039820764525
$ws.Range("Q3").Value = 45.78482962207266
$ws.Range("R3").Value = 412.0634665986539
$ws.Range("S3").Value = 0.02666793997882929
$ws.Range("T3").Value = 0.02666793997882928

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.75560066666667
$ws.Range("H4").Value = 50.266802
$ws.Range("I4").Value = 0.9064438825950116
$ws.Range("J4").Value = 0.9064438825950115
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 82.68775466666666
$ws.Range("N4").Value = 248.063264
$ws.Range("O4").Value = 0.8902831314332034
$ws.Range("P4").Value = 0.8902831314332033
$ws.Range("Q4").Value = 1385.48299721797
$ws.Range("R4").Value = 12469.34697496173
$ws.Range("S4").Value = 0.8069916982651579
$ws.Range("T4").Value = 0.8069916982651577

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.07688299999999999
$ws.Range("H5").Value = 0.230649
$ws.Range("I5").Value = 0.004159213770485276
$ws.Range("J5").Value = 0.004159213770485276
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.457778999999999
$ws.Range("N5").Value = 22.373337
$ws.Range("O5").Value = 0.08029647035915141
$ws.Range("P5").Value = 0.0802964703591514
$ws.Range("Q5").Value = 0.5733764228569999
$ws.Range("R5").Value = 5.160387805712999
$ws.Range("S5").Value = 0.0003339701852391454
$ws.Range("T5").Value = 0.0003339701852391453

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.07688299999999999
$ws.Range("H6").Value = 0.230649
$ws.Range("I6").Value = 0.004159213770485276
$ws.Range("J6").Value = 0.004159213770485276
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.732509
$ws.Range("N6").Value = 8.197526999999999
$ws.Range("O6").Value = 0.02942039820764526
$ws.Range("P6").Value = 0.02942039820764525
$ws.Range("Q6").Value = 0.210083489447
$ws.Range("R6").Value = 1.890751405023
$ws.Range("S6").Value = 0.0001223657253583985
$ws.Range("T6").Value = 0.0001223657253583985

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.07688299999999999
$ws.Range("H7").Value = 0.230649
$ws.Range("I7").Value = 0.004159213770485276
$ws.Range("J7").Value = 0.004159213770485276
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 82.68775466666666
$ws.Range("N7").Value = 248.063264
$ws.Range("O7").Value = 0.8902831314332034
$ws.Range("P7").Value = 0.8902831314332033
$ws.Range("Q7").Value = 6.357282642037332
$ws.Range("R7").Value = 57.215543778336
$ws.Range("S7").Value = 0.003702877859887733
$ws.Range("T7").Value = 0.003702877859887732

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.652500333333333
$ws.Range("H8").Value = 4.957501
$ws.Range("I8").Value = 0.0893969036345032
$ws.Range("J8").Value = 0.08939690363450319
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.457778999999999
$ws.Range("N8").Value = 22.373337
$ws.Range("O8").Value = 0.08029647035915141
$ws.Range("P8").Value = 0.0802964703591514
$ws.Range("Q8").Value = 12.32398228342633
$ws.Range("R8").Value = 110.915840550837
$ws.Range("S8").Value = 0.007178255822887801
$ws.Range("T8").Value = 0.007178255822887799

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.652500333333333
$ws.Range("H9").Value = 4.957501
$ws.Range("I9").Value = 0.0893969036345032
$ws.Range("J9").Value = 0.08939690363450319
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.732509
$ws.Range("N9").Value = 8.197526999999999
$ws.Range("O9").Value = 0.02942039820764526
$ws.Range("P9").Value = 0.02942039820764525
$ws.Range("Q9").Value = 4.515472033336333
$ws.Range("R9").Value = 40.63924830002699
$ws.Range("S9").Value = 0.002630092503457574
$ws.Range("T9").Value = 0.002630092503457573

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.652500333333333
$ws.Range("H10").Value = 4.957501
$ws.Range("I10").Value = 0.0893969036345032
$ws.Range("J10").Value = 0.08939690363450319
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 82.68775466666666
$ws.Range("N10").Value = 248.063264
$ws.Range("O10").Value = 0.8902831314332034
$ws.Range("P10").Value = 0.8902831314332033
$ws.Range("Q10").Value = 136.6415421492515
$ws.Range("R10").Value = 1229.773879343264
$ws.Range("S10").Value = 0.07958855530815782
$ws.Range("T10").Value = 0.07958855530815781
